$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add newly-downloaded metadata values (Satinsky 2015 bac AB / chl Macapa, May 11)
$ws.Range("G24").Value = 3
$ws.Range("L26").Value = 1
$ws.Range("M26").Value = 1

# Update the view: scroll position and selection, matching the edited sheet
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G29:M29").Select()
